# Auto-generated Excel COM-interop script to apply scheduled market-price refresh
# to the Cactuar leve-profit workbook. For each worksheet, the cached price columns
# (H-N) are updated cell-by-cell to match the latest Universalis price snapshot.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2568601.5
$ws.Range("J43").Value = 5665.3335
$ws.Range("L43").Value = 5665.3335
$ws.Range("N43").Value = -5803.3335
$ws.Range("H51").Value = 4789.0835
$ws.Range("I51").Value = 5007.143
$ws.Range("J51").Value = 4483.8
$ws.Range("K51").Value = 5007.143
$ws.Range("L51").Value = 4483.8
$ws.Range("M51").Value = -4523.143
$ws.Range("N51").Value = -5451.8
$ws.Range("H62").Value = 3924.625
$ws.Range("J62").Value = 4279.4
$ws.Range("L62").Value = 4279.4
$ws.Range("N62").Value = -5527.4
$ws.Range("H65").Value = 3924.625
$ws.Range("J65").Value = 4279.4
$ws.Range("L65").Value = 21397
$ws.Range("N65").Value = -27637
$ws.Range("H98").Value = 970.2143
$ws.Range("I98").Value = 798.5833
$ws.Range("J98").Value = 2000
$ws.Range("K98").Value = 798.5833
$ws.Range("L98").Value = 2000
$ws.Range("M98").Value = 699.4167
$ws.Range("N98").Value = -4996
$ws.Range("H106").Value = 12822070
$ws.Range("I106").Value = 13890103
$ws.Range("K106").Value = 13890103
$ws.Range("M106").Value = -13889472
$ws.Range("H107").Value = 799.75
$ws.Range("I107").Value = 731
$ws.Range("J107").Value = 1006
$ws.Range("K107").Value = 731
$ws.Range("L107").Value = 1006
$ws.Range("M107").Value = 1189
$ws.Range("N107").Value = -4846
$ws.Range("H122").Value = 970.2143
$ws.Range("I122").Value = 798.5833
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 2395.7499
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = 54.2501000000002
$ws.Range("N122").Value = -10900
$ws.Range("H132").Value = 13556.103
$ws.Range("I132").Value = 2621.9312
$ws.Range("K132").Value = 7865.7936
$ws.Range("M132").Value = -5335.7936
$ws.Range("H136").Value = 129260
$ws.Range("J136").Value = 129260
$ws.Range("L136").Value = 129260
$ws.Range("N136").Value = -139460
$ws.Range("H137").Value = 11116511
$ws.Range("I137").Value = 2491.8
$ws.Range("J137").Value = 13339315
$ws.Range("K137").Value = 7475.400000000001
$ws.Range("L137").Value = 40017945
$ws.Range("M137").Value = -4925.400000000001
$ws.Range("N137").Value = -40023045
$ws.Range("H138").Value = 2893.1167
$ws.Range("J138").Value = 3337.4043
$ws.Range("L138").Value = 10012.2129
$ws.Range("N138").Value = -20292.2129
$ws.Range("H140").Value = 61870.89
$ws.Range("J140").Value = 60766.125
$ws.Range("L140").Value = 60766.125
$ws.Range("N140").Value = -71126.125

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5368.4424
$ws.Range("I32").Value = 3481.6216
$ws.Range("J32").Value = 10022.6
$ws.Range("K32").Value = 3481.6216
$ws.Range("L32").Value = 10022.6
$ws.Range("M32").Value = -3194.6216
$ws.Range("N32").Value = -10596.6
$ws.Range("H46").Value = 10687
$ws.Range("J46").Value = 12039.182
$ws.Range("L46").Value = 12039.182
$ws.Range("N46").Value = -12677.182
$ws.Range("H122").Value = 3565.375
$ws.Range("I122").Value = 2392.2
$ws.Range("K122").Value = 7176.599999999999
$ws.Range("M122").Value = -4726.599999999999
$ws.Range("H132").Value = 3190.2563
$ws.Range("I132").Value = 1776.3
$ws.Range("J132").Value = 4678.6313
$ws.Range("K132").Value = 5328.9
$ws.Range("L132").Value = 14035.8939
$ws.Range("M132").Value = -2798.9
$ws.Range("N132").Value = -19095.8939

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H56").Value = 15000
$ws.Range("I56").Value = 15000
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 15000
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -14261
$ws.Range("N56").Value = $null
$ws.Range("H97").Value = 8104.5
$ws.Range("I97").Value = 8104.5
$ws.Range("K97").Value = 8104.5
$ws.Range("M97").Value = -7113.5
$ws.Range("H99").Value = 1042938.56
$ws.Range("I99").Value = 1489271.9
$ws.Range("J99").Value = 1494
$ws.Range("K99").Value = 1489271.9
$ws.Range("L99").Value = 1494
$ws.Range("M99").Value = -1487773.9
$ws.Range("N99").Value = -4490
$ws.Range("H134").Value = 6330.75
$ws.Range("I134").Value = 3666.5
$ws.Range("J134").Value = 8995
$ws.Range("K134").Value = 10999.5
$ws.Range("L134").Value = 26985
$ws.Range("M134").Value = -8464.5
$ws.Range("N134").Value = -32055

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").Value = $null
$ws.Range("H132").Value = 55580284
$ws.Range("I132").Value = 83352590
$ws.Range("K132").Value = 250057770
$ws.Range("M132").Value = -250055240
$ws.Range("H141").Value = 106554.9
$ws.Range("J141").Value = 106554.9
$ws.Range("L141").Value = 106554.9
$ws.Range("N141").Value = -116914.9

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 39.5
$ws.Range("I12").Value = 39.5
$ws.Range("K12").Value = 118.5
$ws.Range("M12").Value = 54.5
$ws.Range("H50").Value = 628.4286
$ws.Range("I50").Value = 459.8
$ws.Range("J50").Value = 1050
$ws.Range("K50").Value = 1379.4
$ws.Range("L50").Value = 3150
$ws.Range("M50").Value = -898.4000000000001
$ws.Range("N50").Value = -4112
$ws.Range("H53").Value = 628.4286
$ws.Range("I53").Value = 459.8
$ws.Range("J53").Value = 1050
$ws.Range("K53").Value = 1379.4
$ws.Range("L53").Value = 3150
$ws.Range("M53").Value = -898.4000000000001
$ws.Range("N53").Value = -4112
$ws.Range("H75").Value = 5109.923
$ws.Range("J75").Value = 5611.727
$ws.Range("L75").Value = 16835.181
$ws.Range("N75").Value = -18831.181
$ws.Range("H78").Value = 5109.923
$ws.Range("J78").Value = 5611.727
$ws.Range("L78").Value = 50505.543
$ws.Range("N78").Value = -60489.543
$ws.Range("H122").Value = 6452311.5
$ws.Range("I122").Value = 16129480
$ws.Range("J122").Value = 865.6667
$ws.Range("K122").Value = 145165320
$ws.Range("L122").Value = 7791.0003
$ws.Range("M122").Value = -145162870
$ws.Range("N122").Value = -12691.0003
$ws.Range("H131").Value = 8477158
$ws.Range("J131").Value = 6947123
$ws.Range("L131").Value = 20841369
$ws.Range("N131").Value = -20851449
$ws.Range("H134").Value = 23618.545
$ws.Range("I134").Value = 19838.334
$ws.Range("K134").Value = 59515.00199999999
$ws.Range("M134").Value = -54445.00199999999

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 292.16666
$ws.Range("I97").Value = 259.25
$ws.Range("K97").Value = 259.25
$ws.Range("M97").Value = 236.75
$ws.Range("H107").Value = 2801803
$ws.Range("I107").Value = 5291576
$ws.Range("K107").Value = 5291576
$ws.Range("M107").Value = -5289656
$ws.Range("H113").Value = 1762.5
$ws.Range("J113").Value = 1762.5
$ws.Range("L113").Value = 1762.5
$ws.Range("N113").Value = -6102.5
$ws.Range("H122").Value = 410406.53
$ws.Range("I122").Value = 552664.1
$ws.Range("J122").Value = 3956.1428
$ws.Range("K122").Value = 1657992.3
$ws.Range("L122").Value = 11868.4284
$ws.Range("M122").Value = -1655542.3
$ws.Range("N122").Value = -16768.4284
$ws.Range("H126").Value = 4793.5654
$ws.Range("I126").Value = 2782.7334
$ws.Range("J126").Value = 8563.875
$ws.Range("K126").Value = 8348.200199999999
$ws.Range("L126").Value = 25691.625
$ws.Range("M126").Value = -5878.200199999999
$ws.Range("N126").Value = -30631.625
$ws.Range("H132").Value = 3232.6572
$ws.Range("I132").Value = 2609.2
$ws.Range("K132").Value = 7827.599999999999
$ws.Range("M132").Value = -5297.599999999999
$ws.Range("H140").Value = 78779.5
$ws.Range("J140").Value = 78779.5
$ws.Range("L140").Value = 78779.5
$ws.Range("N140").Value = -89139.5

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H139").Value = 89715
$ws.Range("J139").Value = 89715
$ws.Range("L139").Value = 89715
$ws.Range("N139").Value = -99995

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 39999.332
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 39999.332
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 39999.332
$ws.Range("M62").Value = $null
$ws.Range("N62").Value = -41247.332
$ws.Range("H65").Value = 39999.332
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 39999.332
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 199996.66
$ws.Range("M65").Value = $null
$ws.Range("N65").Value = -206236.66
$ws.Range("H102").Value = 93332.336
$ws.Range("J102").Value = 93332.336
$ws.Range("L102").Value = 93332.336
$ws.Range("N102").Value = -99822.336
$ws.Range("H125").Value = 80000
$ws.Range("J125").Value = 80000
$ws.Range("L125").Value = 80000
$ws.Range("N125").Value = -89840
$ws.Range("H126").Value = 3302.75
$ws.Range("I126").Value = 3541.5715
$ws.Range("J126").Value = 2968.4
$ws.Range("K126").Value = 10624.7145
$ws.Range("L126").Value = 8905.200000000001
$ws.Range("M126").Value = -8154.7145
$ws.Range("N126").Value = -13845.2
$ws.Range("H136").Value = 9547.209000000001
$ws.Range("J136").Value = 11860.465
$ws.Range("L136").Value = 35581.395
$ws.Range("N136").Value = -40681.395
$ws.Range("H141").Value = 49837
$ws.Range("J141").Value = 49813.715
$ws.Range("L141").Value = 49813.715
$ws.Range("N141").Value = -60173.715
